$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 15.29788533997392

$ws.Range("I4").Value = 19.12364428961322

$ws.Range("I5").Value = 15.20442965828532

$ws.Range("H6").Value = 15292167840.66753
$ws.Range("J6").Value = 15275111475.24988

$ws.Range("H7").Value = 15292167840.66753
$ws.Range("J7").Value = 15275111475.24988

$ws.Range("H8").Value = 15292160198.18902
$ws.Range("I8").Value = 14.48403123646112
$ws.Range("J8").Value = 15275221470.16526

$ws.Range("H9").Value = 15292160198.18902
$ws.Range("I9").Value = 14.48403123646112
$ws.Range("J9").Value = 15275221470.16526
